$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 606.9545000000001
$ws.Range("I28").Value = 585.3889
$ws.Range("J28").Value = 704
$ws.Range("K28").Value = 585.3889
$ws.Range("L28").Value = 704
$ws.Range("M28").Value = -100.3889
$ws.Range("N28").Value = -1674

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H40").Value = 1692.3462
$ws.Range("I40").Value = 1005.6923
$ws.Range("J40").Value = 5125.615
$ws.Range("K40").Value = 1005.6923
$ws.Range("L40").Value = 5125.615
$ws.Range("M40").Value = -830.6923
$ws.Range("N40").Value = -5475.615

$ws.Range("H43").Value = 1774.875
$ws.Range("I43").Value = 2366.6667
$ws.Range("J43").Value = 1419.8
$ws.Range("K43").Value = 2366.6667
$ws.Range("L43").Value = 1419.8
$ws.Range("M43").Value = -2297.6667
$ws.Range("N43").Value = -1557.8

$ws.Range("H55").Value = 144.71428
$ws.Range("I55").Value = 148.83333
$ws.Range("J55").Value = 120
$ws.Range("K55").Value = 148.83333
$ws.Range("L55").Value = 120
$ws.Range("M55").Value = 65.16667000000001
$ws.Range("N55").Value = -548

$ws.Range("H75").Value = 30154.166
$ws.Range("J75").Value = 30154.166
$ws.Range("L75").Value = 30154.166
$ws.Range("N75").Value = -32026.166

$ws.Range("H78").Value = 30154.166
$ws.Range("J78").Value = 30154.166
$ws.Range("L78").Value = 90462.49800000001
$ws.Range("N78").Value = -99822.49800000001

$ws.Range("H125").Value = 1536
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1536
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 13824
$ws.Range("N125").Value = -18744
$ws.Range("M125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1086
$ws.Range("I2").Value = 880.875
$ws.Range("J2").Value = 1906.5
$ws.Range("K2").Value = 880.875
$ws.Range("L2").Value = 1906.5
$ws.Range("M2").Value = -767.875
$ws.Range("N2").Value = -2132.5

$ws.Range("H32").Value = 1211897.4
$ws.Range("I32").Value = 1304438.1
$ws.Range("J32").Value = 27376.6
$ws.Range("K32").Value = 1304438.1
$ws.Range("L32").Value = 27376.6
$ws.Range("M32").Value = -1304151.1
$ws.Range("N32").Value = -27950.6

$ws.Range("H45").Value = 1097
$ws.Range("I45").Value = 1032.5454
$ws.Range("J45").Value = 1333.3334
$ws.Range("K45").Value = 1032.5454
$ws.Range("L45").Value = 1333.3334
$ws.Range("M45").Value = -655.5454
$ws.Range("N45").Value = -2087.3334

$ws.Range("H61").Value = 394334.16
$ws.Range("I61").Value = 334671.34
$ws.Range("J61").Value = 479566.75
$ws.Range("K61").Value = 334671.34
$ws.Range("L61").Value = 479566.75
$ws.Range("M61").Value = -334459.34
$ws.Range("N61").Value = -479990.75

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 270279.28
$ws.Range("I74").Value = 371571.72
$ws.Range("J74").Value = 74929.57000000001
$ws.Range("K74").Value = 371571.72
$ws.Range("L74").Value = 74929.57000000001
$ws.Range("M74").Value = -370697.72
$ws.Range("N74").Value = -76677.57000000001

$ws.Range("H77").Value = 270279.28
$ws.Range("I77").Value = 371571.72
$ws.Range("J77").Value = 74929.57000000001
$ws.Range("K77").Value = 1857858.6
$ws.Range("L77").Value = 374647.85
$ws.Range("M77").Value = -1853490.6
$ws.Range("N77").Value = -383383.85

$ws.Range("H110").Value = 2425.4546
$ws.Range("I110").Value = 2368
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 2368
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = -323
$ws.Range("N110").Value = -7090

$ws.Range("H116").Value = 1086
$ws.Range("I116").Value = 880.875
$ws.Range("J116").Value = 1906.5
$ws.Range("K116").Value = 880.875
$ws.Range("L116").Value = 1906.5
$ws.Range("M116").Value = 1413.125
$ws.Range("N116").Value = -6494.5

$ws.Range("H132").Value = 30384.893
$ws.Range("I132").Value = 45106.332
$ws.Range("K132").Value = 135318.996
$ws.Range("M132").Value = -132788.996

$ws.Range("H136").Value = 394334.16
$ws.Range("I136").Value = 334671.34
$ws.Range("J136").Value = 479566.75
$ws.Range("K136").Value = 1004014.02
$ws.Range("L136").Value = 1438700.25
$ws.Range("M136").Value = -1001464.02
$ws.Range("N136").Value = -1443800.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1086
$ws.Range("I3").Value = 880.875
$ws.Range("J3").Value = 1906.5
$ws.Range("K3").Value = 880.875
$ws.Range("L3").Value = 1906.5
$ws.Range("M3").Value = -766.875
$ws.Range("N3").Value = -2134.5

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H134").Value = 2014.8918
$ws.Range("I134").Value = 1196.05
$ws.Range("J134").Value = 2978.2354
$ws.Range("K134").Value = 3588.15
$ws.Range("L134").Value = 8934.706200000001
$ws.Range("M134").Value = -1053.15
$ws.Range("N134").Value = -14004.7062

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2112.4792
$ws.Range("I31").Value = 1109.975
$ws.Range("K31").Value = 1109.975
$ws.Range("M31").Value = -814.9749999999999

$ws.Range("H34").Value = 2112.4792
$ws.Range("I34").Value = 1109.975
$ws.Range("K34").Value = 1109.975
$ws.Range("M34").Value = -907.9749999999999

$ws.Range("H105").Value = 910.625
$ws.Range("I105").Value = 880.9524
$ws.Range("J105").Value = 1118.3334
$ws.Range("K105").Value = 880.9524
$ws.Range("L105").Value = 1118.3334
$ws.Range("M105").Value = 866.0476
$ws.Range("N105").Value = -4612.3334

$ws.Range("H107").Value = 430.63416
$ws.Range("I107").Value = 370.96667
$ws.Range("J107").Value = 593.36365
$ws.Range("K107").Value = 370.96667
$ws.Range("L107").Value = 593.36365
$ws.Range("M107").Value = 1549.03333
$ws.Range("N107").Value = -4433.36365

$ws.Range("H134").Value = 2155.652
$ws.Range("I134").Value = 1256.25
$ws.Range("J134").Value = 2635.3333
$ws.Range("K134").Value = 3768.75
$ws.Range("L134").Value = 7905.999899999999
$ws.Range("M134").Value = -1233.75
$ws.Range("N134").Value = -12975.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1248.92
$ws.Range("I5").Value = 463.45456
$ws.Range("K5").Value = 1390.36368
$ws.Range("M5").Value = -1278.36368

$ws.Range("H12").Value = 44.25
$ws.Range("J12").Value = 50.57143
$ws.Range("L12").Value = 151.71429
$ws.Range("N12").Value = -497.71429

$ws.Range("H74").Value = 4999.3335
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 6999
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 20997
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -23119

$ws.Range("H77").Value = 4999.3335
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 6999
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 62991
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -73599

$ws.Range("H122").Value = 23810078
$ws.Range("J122").Value = 835.4
$ws.Range("L122").Value = 7518.599999999999
$ws.Range("N122").Value = -12418.6

$ws.Range("H135").Value = 1248.92
$ws.Range("I135").Value = 463.45456
$ws.Range("K135").Value = 4171.09104
$ws.Range("M135").Value = -1636.09104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3772.5227
$ws.Range("I80").Value = 4237.4136
$ws.Range("K80").Value = 4237.4136
$ws.Range("M80").Value = -3239.4136

$ws.Range("H83").Value = 3772.5227
$ws.Range("I83").Value = 4237.4136
$ws.Range("K83").Value = 21187.068
$ws.Range("M83").Value = -16195.068

$ws.Range("H102").Value = 3771.158
$ws.Range("I102").Value = 2919
$ws.Range("J102").Value = 4390.909
$ws.Range("K102").Value = 2919
$ws.Range("L102").Value = 4390.909
$ws.Range("M102").Value = -1297
$ws.Range("N102").Value = -7634.909

$ws.Range("H113").Value = 880.0909
$ws.Range("I113").Value = 553.2222
$ws.Range("J113").Value = 2351
$ws.Range("K113").Value = 553.2222
$ws.Range("L113").Value = 2351
$ws.Range("M113").Value = 1616.7778
$ws.Range("N113").Value = -6691

$ws.Range("H132").Value = 2811.1667
$ws.Range("I132").Value = 2757.0732
$ws.Range("J132").Value = 2871.1082
$ws.Range("K132").Value = 8271.2196
$ws.Range("L132").Value = 8613.3246
$ws.Range("M132").Value = -5741.2196
$ws.Range("N132").Value = -13673.3246

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19148.334
$ws.Range("I40").Value = 22378
$ws.Range("K40").Value = 22378
$ws.Range("M40").Value = -22242

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 9812.5
$ws.Range("I132").Value = 2572.3125
$ws.Range("J132").Value = 19466.084
$ws.Range("K132").Value = 7716.9375
$ws.Range("L132").Value = 58398.25199999999
$ws.Range("M132").Value = -5186.9375
$ws.Range("N132").Value = -63458.25199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 685.7143
$ws.Range("I107").Value = 540
$ws.Range("K107").Value = 1620
$ws.Range("M107").Value = 300
